# Update the "想去人数" (number of people interested) figures that changed
# between the two gh-pages data generation runs.
#
# Sheet "展览"   (Exhibitions)
# Sheet "演出"   (Performances)
# Sheet "本地生活" (Local life) -- unchanged
# Sheet "全部类型" (All types)

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 9320
$wsExhibition.Range("F4").Value = 17
$wsExhibition.Range("F5").Value = 503
$wsExhibition.Range("F6").Value = 459

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 1

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 9320
$wsAll.Range("F4").Value = 17
$wsAll.Range("F5").Value = 503
$wsAll.Range("F6").Value = 1
$wsAll.Range("F7").Value = 459
